$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values
$ws.Range("C2").Value = 11
$ws.Range("B3").Value = 5.5
$ws.Range("C3").Value = 9.5

# Update the selected cell (active cell / selection) to C2
$ws.Range("C2").Select()
